# Add results 001 - 004 (5 new log rows: 22-26) to the "logs" worksheet

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("logs")

$bNameFeatures = "12 features: %ascii-adp, %digit-adp, digit-adp/ascii-adp, %keyword-name, %keyword-address, %keyword-phone, b#max-digit-skip-all-punctuation >= 7, bfirst-character-digit, bfirst-character-ascii, blast-character-digit, blast-character-ascii, b#ascii >= 6"
$cAddressFeatures = "11 features: %ascii-adp, %digit-adp, digit-adp/ascii-adp, %keyword-name, %keyword-address, %keyword-phone, b#max-digit-skip-all-punctuation >= 7, bfirst-character-digit, bfirst-character-ascii, blast-character-digit, blast-character-ascii"
$dPhoneFeatures = "12 features: %ascii-adp, %digit-adp, digit-adp/ascii-adp, %keyword-name, %keyword-address, %keyword-phone, b#max-digit-skip-all-punctuation >= 7, bfirst-character-digit, bfirst-character-ascii, blast-character-digit, blast-character-ascii, b#digit >= 7"
$modelType = "Neural-Network"
$model = "2 layers: [10-Sigmoid, 2-Softmax], learning_rate: 0.01, learning_rule: adagrad, n_iterator: 3000"

$rows = @(
    @{ Time = "20160427_084027"; ClassifyAcc = 0.891089108910891; SegmentAcc = 0.42 },
    @{ Time = "20160427_092248"; ClassifyAcc = 0.897689768976898; SegmentAcc = 0.43 },
    @{ Time = "20160427_100536"; ClassifyAcc = 0.900990099009901; SegmentAcc = 0.45 },
    @{ Time = "20160427_104748"; ClassifyAcc = 0.900990099009901; SegmentAcc = 0.44 },
    @{ Time = "20160427_113113"; ClassifyAcc = 0.907590759075908; SegmentAcc = 0.45 }
)

$r = 22
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = $row.Time
    $ws.Cells.Item($r, 2).Value = $bNameFeatures
    $ws.Cells.Item($r, 3).Value = $cAddressFeatures
    $ws.Cells.Item($r, 4).Value = $dPhoneFeatures
    $ws.Cells.Item($r, 5).Value = $modelType
    $ws.Cells.Item($r, 6).Value = $model
    $ws.Cells.Item($r, 7).Value = $modelType
    $ws.Cells.Item($r, 8).Value = $model
    $ws.Cells.Item($r, 9).Value = $modelType
    $ws.Cells.Item($r, 10).Value = $model
    $ws.Cells.Item($r, 11).Value = $row.ClassifyAcc
    $ws.Cells.Item($r, 12).Value = $row.SegmentAcc
    $r = $r + 1
}
